# Generate Report for Handback
# Update the timestamp values recorded on the handback-status report.
# These cells hold plain text timestamps (formatted to look like
# dates), so we explicitly force text so Excel doesn't reinterpret
# them as date serials.

$wb = $excel.ActiveWorkbook

# "Overview" sheet: Latest HO Xliff Generate Date
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-08-30 21:13:37"

# "zh-cn" sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-08-30 21:13:31"
$wsZhCn.Range("K2").Value = "2016-08-30 21:13:50"

# "de-de" sheet: Correspond Handoff Datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-08-30 21:13:57"
